$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of (row, newValue) pairs, matching the order in which the
# new shared strings should be appended to the workbook's string table.
$updates = @(
    @(6,  "a1-reffacility"),
    @(7,  "a1-refname"),
    @(8,  "a1-id"),
    @(9,  "a1-name"),
    @(10, "a1-sex"),
    @(11, "a1-enroldate"),
    @(12, "a4-a4_b_2"),
    @(13, "a4-a4_b_6"),
    @(14, "a4-a4_b_3"),
    @(15, "a4-a4_b_1a"),
    @(16, "a4-a4_b_1"),
    @(17, "a1-hvisitdate"),
    @(18, "a1-found"),
    @(19, "a1-n2_1"),
    @(20, "a1-n2_2"),
    @(21, "n2t-n3_1a"),
    @(22, "n2-t_n3_1"),
    @(23, "n2t-n3_2a"),
    @(24, "n2t-n3_2"),
    @(25, "n4-n2_3"),
    @(26, "n4-n4_4"),
    @(27, "n4-n4_4o"),
    @(28, "n4-n4_1"),
    @(29, "n4-n4_1o"),
    @(30, "n4-n4_2a"),
    @(31, "n4-n4_2"),
    @(32, "n4-n4_3"),
    @(33, "n4-n4_5"),
    @(5,  "SubmitterID")
)

foreach ($pair in $updates) {
    $row = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($row, 1).Value = $val
}

$null = $ws.Range("E7").Select()
